# Rename the worksheet from "Sheet1" to "listOfCities"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "listOfCities"

# Clear the stored A5 selection, resetting it back to the default
# top-left cell (A1).
$ws.Range("A1").Select()
